$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Set column F (variance) values for rows 2 through 52 to 2,
# with the exception of rows 11 and 12 which take special values.
for ($row = 2; $row -le 52; $row++) {
    if ($row -eq 11) {
        $ws.Cells.Item($row, 6).Value = 50
    } elseif ($row -eq 12) {
        $ws.Cells.Item($row, 6).Value = 100
    } else {
        $ws.Cells.Item($row, 6).Value = 2
    }
}

# Reflect the final cell selection left behind in the saved sheet view.
$ws.Range("G9").Select()
